# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header -> "Weekly_PO_Qty" on "Weekly Quantity" sheet
# 2. Rename "Requested quantity" header -> "Monthly_PO_Qty" on "Monthly Trend" sheet
# 3. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data

$wb = $excel.ActiveWorkbook

# --- 1 & 2: relabel the "Requested quantity" header cells -----------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 3: add the "PO Forecast" sheet after the last existing sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header formatting (bold / bordered / centered style) from the
# "Weekly Quantity" header row, and the date-column formatting from its
# first data cell, so the new sheet reuses the workbook's existing styles.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A25").PasteSpecial(-4122)

$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

$forecastData = @(
    @(45347.99999999999, 3, 0.2694875298912611, 5.562509485642984),
    @(45396.99999999999, 3, 0.172702610333869, 5.574884948193354),
    @(45410.99999999999, 3, 0.1645390559312308, 5.716446041579286),
    @(45417.99999999999, 3, 0.38886146198068, 5.71060761360079),
    @(45424.99999999999, 3, 0.4895285253757932, 5.761034819544324),
    @(45438.99999999999, 3, 0.4136874077571332, 5.636279921486546),
    @(45487.99999999999, 3, 0.1158866493743832, 5.518297605951786),
    @(45501.99999999999, 3, 0.3133322065165772, 5.837923231456047),
    @(45508.99999999999, 3, 0.3727395695877033, 5.710113326182914),
    @(45515.99999999999, 3, 0.1328794448939764, 5.742927235209355),
    @(45522.99999999999, 3, 0.4364889734049733, 5.697223472554938),
    @(45529.99999999999, 3, 0.2093725055266555, 5.740329141289823),
    @(45543.99999999999, 3, 0.2783334395222958, 5.563223768283376),
    @(45550.99999999999, 3, 0.2082365715622867, 5.608512892404942),
    @(45571.99999999999, 3, 0.281153398961461, 5.54677762619685),
    @(45578.99999999999, 3, 0.5436025362293007, 5.684070877943914),
    @(45585.99999999999, 3, 0.2289732462401197, 5.804838491676238),
    @(45592.99999999999, 3, 0.3158940551873764, 5.802681459224644),
    @(45599.99999999999, 3, 0.511437387783285, 5.627238316411318),
    @(45606.99999999999, 3, 0.4758554760535734, 5.88874157761177),
    @(45613.99999999999, 3, 0.3002359049040087, 5.776942773405858),
    @(45620.99999999999, 3, 0.5471116525303314, 5.891853751857901),
    @(45627.99999999999, 3, 0.3203319661845824, 5.711481125951017),
    @(45634.99999999999, 3, 0.3633552452880243, 5.607250513117183)
)

for ($i = 0; $i -lt $forecastData.Length; $i++) {
    $row = $i + 2
    $wsForecast.Range("A$row").Value = $forecastData[$i][0]
    $wsForecast.Range("B$row").Value = $forecastData[$i][1]
    $wsForecast.Range("C$row").Value = $forecastData[$i][2]
    $wsForecast.Range("D$row").Value = $forecastData[$i][3]
}
